# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
# This updates column G ("K") values for rows 2-38 (row 3 stays unchanged at 0).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2  = 1
    4  = 6
    5  = 9
    6  = 3
    7  = 2
    8  = 5
    9  = 4
    10 = 9
    11 = 2
    12 = 6
    13 = 5
    14 = 1
    15 = 4
    16 = 8
    17 = 5
    18 = 7
    19 = 6
    20 = 2
    21 = 4
    22 = 7
    23 = 4
    24 = 5
    25 = 6
    26 = 5
    27 = 9
    28 = 5
    29 = 4
    30 = 9
    31 = 6
    32 = 6
    33 = 2
    34 = 5
    35 = 3
    36 = 6
    37 = 3
    38 = 4
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}
